$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current "ASIN" column (B) to hold the
# new "Week_Start_Date" field. This shifts ASIN..is_holiday_week one
# column to the right (B:I -> C:J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Per-row week start dates (keep them as literal text, like the rest of
# the sheet's text columns, not auto-converted Excel date serials).
$ws.Range("B2").Value = "'2025-01-05"
$ws.Range("B3").Value = "'2025-01-12"
$ws.Range("B4").Value = "'2025-01-19"
$ws.Range("B5").Value = "'2025-01-26"
$ws.Range("B6").Value = "'2025-02-02"
$ws.Range("B7").Value = "'2025-02-09"
$ws.Range("B8").Value = "'2025-02-16"
$ws.Range("B9").Value = "'2025-02-23"
$ws.Range("B10").Value = "'2025-03-02"
$ws.Range("B11").Value = "'2025-03-09"
$ws.Range("B12").Value = "'2025-03-16"
$ws.Range("B13").Value = "'2025-03-23"
$ws.Range("B14").Value = "'2025-03-30"
$ws.Range("B15").Value = "'2025-04-06"
$ws.Range("B16").Value = "'2025-04-13"
$ws.Range("B17").Value = "'2025-04-20"

# The Week labels for the first nine weeks drop their leading zero
# (W01 -> W1 ... W09 -> W9); W10-W16 were already unpadded.
$ws.Range("A2").Value = "W1"
$ws.Range("A3").Value = "W2"
$ws.Range("A4").Value = "W3"
$ws.Range("A5").Value = "W4"
$ws.Range("A6").Value = "W5"
$ws.Range("A7").Value = "W6"
$ws.Range("A8").Value = "W7"
$ws.Range("A9").Value = "W8"
$ws.Range("A10").Value = "W9"

# is_holiday_week (now column J) becomes a proper boolean column.
$ws.Range("J2").Value = $false
$ws.Range("J3").Value = $false
$ws.Range("J4").Value = $false
$ws.Range("J5").Value = $false
$ws.Range("J6").Value = $false
$ws.Range("J7").Value = $false
$ws.Range("J8").Value = $false
$ws.Range("J9").Value = $false
$ws.Range("J10").Value = $false
$ws.Range("J11").Value = $false
$ws.Range("J12").Value = $false
$ws.Range("J13").Value = $false
$ws.Range("J14").Value = $false
$ws.Range("J15").Value = $false
$ws.Range("J16").Value = $false
$ws.Range("J17").Value = $false
